$wb = $excel.ActiveWorkbook

# Record "1475" (dated 2025-11-10, Merkez İlçe, job 3B,
# HAVVA NİLGÜN KIYMAÇ (K.Mühendisi), NURHAN ARSLAN (Tekniker)) was deleted.
# It appears on the "Kayitlar" sheet (row 601) and is mirrored on the
# "Merkez İlçe" sheet (row 59). Delete the whole row in both places so
# every subsequent row shifts up by one.

$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
$wsKayitlar.Rows.Item(601).Delete()

$wsMerkez = $wb.Worksheets.Item("Merkez İlçe")
$wsMerkez.Rows.Item(59).Delete()
